$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "28.134.11"
Set-TextValue 2 5 "  +0.01%  "

Set-TextValue 3 4 "1.882.96"
Set-TextValue 3 5 "  -0.78%  "

Set-TextValue 4 5 "  +0.32%  "

Set-TextValue 5 5 "  -0.28%  "

Set-TextValue 6 5 "  +0.30%  "

Set-TextValue 7 4 "0.5044"
Set-TextValue 7 5 "  +0.20%  "

Set-TextValue 8 5 "  -1.94%  "

Set-TextValue 9 4 "0.08554"
Set-TextValue 9 5 "  -7.38%  "

Set-TextValue 10 4 "1.115"
Set-TextValue 10 5 "  -1.33%  "

Set-TextValue 11 4 "41.80"
Set-TextValue 11 5 "  -0.03%  "

Set-TextValue 12 5 "  -2.06%  "

Set-TextValue 13 2 "WrappedEther"
Set-TextValue 13 3 "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue 13 4 "1.875.43"
Set-TextValue 13 5 "  -1.57%  "

Set-TextValue 14 2 "Solana"
Set-TextValue 14 3 "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue 14 4 "20.59"
Set-TextValue 14 5 "  -1.07%  "

Set-TextValue 15 4 "7.201"
Set-TextValue 15 5 "  -1.31%  "

Set-TextValue 16 5 "  +0.32%  "

Set-TextValue 17 4 "0.00001099"
Set-TextValue 17 5 "  -1.07%  "

Set-TextValue 18 4 "91.06"
Set-TextValue 18 5 "  -1.46%  "

Set-TextValue 19 4 "0.06662"
Set-TextValue 19 5 "  +0.49%  "

Set-TextValue 20 4 "18.03"
Set-TextValue 20 5 "  +0.84%  "

Set-TextValue 21 4 "1.003"
Set-TextValue 21 5 "  +0.35%  "

Set-TextValue 22 4 "6.091"
Set-TextValue 22 5 "  -2.00%  "

Set-TextValue 23 4 "28.169.58"
Set-TextValue 23 5 "  -0.06%  "

Set-TextValue 24 4 "11.18"
Set-TextValue 24 5 "  -2.58%  "

Set-TextValue 25 4 "2.269"
Set-TextValue 25 5 "  -2.13%  "

Set-TextValue 26 4 "2.579"
Set-TextValue 26 5 "  +1.04%  "

Set-TextValue 27 4 "2.096.23"
Set-TextValue 27 5 "  -1.37%  "

Set-TextValue 28 4 "20.66"
Set-TextValue 28 5 "  -1.05%  "

Set-TextValue 29 4 "156.21"
Set-TextValue 29 5 "  -1.48%  "

Set-TextValue 30 4 "126.42"
Set-TextValue 30 5 "  -0.43%  "

Set-TextValue 31 4 "0.1051"
Set-TextValue 31 5 "  -1.12%  "

Set-TextValue 32 4 "1.048"
Set-TextValue 32 5 "  -3.56%  "

Set-TextValue 33 4 "5.626"
Set-TextValue 33 5 "  +0.15%  "

Set-TextValue 34 4 "3.605"
Set-TextValue 34 5 "  -0.28%  "

Set-TextValue 35 4 "9.674"
Set-TextValue 35 5 "  +0.67%  "

Set-TextValue 36 4 "0.02450"
Set-TextValue 36 5 "  +2.17%  "

Set-TextValue 37 4 "0.06525"
Set-TextValue 37 5 "  -1.32%  "

Set-TextValue 38 4 "1.232"
Set-TextValue 38 5 "  +0.13%  "

Set-TextValue 39 4 "0.2174"
Set-TextValue 39 5 "  -1.68%  "

Set-TextValue 40 4 "1.238"
Set-TextValue 40 5 "  -8.34%  "

Set-TextValue 41 4 "0.6388"
Set-TextValue 41 5 "  -1.16%  "

Set-TextValue 42 4 "11.37"
Set-TextValue 42 5 "  -0.37%  "

Set-TextValue 43 4 "4.885"
Set-TextValue 43 5 "  -1.79%  "

Set-TextValue 44 4 "0.6081"
Set-TextValue 44 5 "  -0.18%  "

Set-TextValue 45 4 "13.06"
Set-TextValue 45 5 "  -2.59%  "

Set-TextValue 46 4 "1.300"
Set-TextValue 46 5 "  -0.67%  "

Set-TextValue 47 4 "3.680"
Set-TextValue 47 5 "  -0.28%  "

Set-TextValue 48 5 "  -0.31%  "

Set-TextValue 49 4 "1.218"
Set-TextValue 49 5 "  +1.91%  "

Set-TextValue 50 4 "120.67"
Set-TextValue 50 5 "  -1.17%  "

Set-TextValue 51 4 "80.60"
Set-TextValue 51 5 "  +1.95%  "
